$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.893.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.809.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4603"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3702"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8738"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.358"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.528"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07041"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008695"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.896.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.329"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.052.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.902"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.140"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.303"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08898"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7525"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.156"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.441"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.905"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.099"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01971"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05234"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.427"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.930"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5306"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.181"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1664"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.509"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4972"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.672"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("E51").Value = "  -1.65%  "
